# Reception Perception Rookies - add "RP 2021-2025 Rank" / "RP Definition" /
# "RP Quick Note" data (columns E/F/G) for the remaining rookies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Jalen Royals (row 12) ---
$ws.Range("E12").Value = 61
$ws.Range("F12").Value = "Late Round 3  / Early Round 4"
$ws.Range("G12").Value = "Shows flashes of good routes and YAC skills but needs development and early role-catering.`t"
$ws.Range("E12:G12").WrapText = $true

# --- Jack Bech (row 13) ---
$ws.Range("F13").Value = "Priority Rd 2"
$ws.Range("G13").Value = "Love this guy. Beats man and zone. Dominant underneath. Catches everything, breaks tackles.`t"
$ws.Range("F13:G13").WrapText = $true

# --- Tory Horton (row 14) ---
$ws.Range("E14").Value = 37
$ws.Range("F14").Value = "Priority Rd 2"
$ws.Range("G14").Value = "Impressive route runner who projects well as a flanker. Attacks the ball and wins on in-breakers.`t"
$ws.Range("E14:G14").WrapText = $true

# --- Xavier Restrepo (row 15) ---
$ws.Range("E15").Value = 73
$ws.Range("F15").Value = "Late Round 3  / Early Round 4"
$ws.Range("G15").Value = "I can easily see him being a useful popgun slot receiver in the league but limited player.`t"
$ws.Range("E15:G15").WrapText = $true

# --- Savion Williams (row 17) ---
$ws.Range("E17").Value = 87
$ws.Range("F17").Value = "The Rest"
$ws.Range("G17").Value = "He's ranked this low as a pure wide receiver but I think he can be more useful in a gadget role.`t"
$ws.Range("E17:G17").WrapText = $true
$ws.Rows(17).RowHeight = 87

# --- Tez Johnson (row 18) ---
$ws.Range("E18").Value = 75
$ws.Range("F18").Value = "Late Round 3  / Early Round 4"
$ws.Range("G18").Value = "Solid zone-beater but struggles to win against man coverage given lack of size and speed.`t"
$ws.Range("E18:G18").WrapText = $true

# --- Pat Bryant (row 20) ---
$ws.Range("E20").Value = 78
$ws.Range("F20").Value = "Late Round 3  / Early Round 4"
$ws.Range("G20").Value = "I can see him developing into a reliable slot, flanker rotation guy at some point. Has some game.`t"
$ws.Range("E20:G20").WrapText = $true
$ws.Rows(20).RowHeight = 87

# --- Nick Nash (row 21) ---
$ws.Range("E21").Value = 83
$ws.Range("F21").Value = "Late Round 3  / Early Round 4"
$ws.Range("G21").Value = "Will be a 26-year-old rookie but shows some ability to translate into a big slot role.`t"
$ws.Range("F21:G21").WrapText = $true
$ws.Rows(21).RowHeight = 72.5
